# Scheduled runner update: refresh cached FFXIV market-board price/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3918.6667
$ws.Range("I98").Value = 3575.3845
$ws.Range("K98").Value = 3575.3845
$ws.Range("M98").Value = -2077.3845
$ws.Range("H113").Value = 111200970
$ws.Range("I113").Value = 115243.57
$ws.Range("J113").Value = 500001000
$ws.Range("K113").Value = 115243.57
$ws.Range("L113").Value = 500001000
$ws.Range("M113").Value = -111989.57
$ws.Range("N113").Value = -500007508
$ws.Range("H121").Value = 1127.2
$ws.Range("I121").Value = 575
$ws.Range("J121").Value = 1495.3334
$ws.Range("K121").Value = 1725
$ws.Range("L121").Value = 4486.0002
$ws.Range("M121").Value = 22
$ws.Range("N121").Value = -7980.0002
$ws.Range("H122").Value = 3918.6667
$ws.Range("I122").Value = 3575.3845
$ws.Range("K122").Value = 10726.1535
$ws.Range("M122").Value = -8276.1535
$ws.Range("H132").Value = 6174083
$ws.Range("I132").Value = 7247548.5
$ws.Range("K132").Value = 21742645.5
$ws.Range("M132").Value = -21740115.5
$ws.Range("H137").Value = 918.38635
$ws.Range("I137").Value = 710.03705
$ws.Range("K137").Value = 2130.11115
$ws.Range("M137").Value = 419.8888499999998
$ws.Range("H138").Value = 1610.3939
$ws.Range("I138").Value = 1323.8966
$ws.Range("J138").Value = 2015.683
$ws.Range("K138").Value = 3971.6898
$ws.Range("L138").Value = 6047.049
$ws.Range("M138").Value = 1168.3102
$ws.Range("N138").Value = -16327.049
$ws.Range("H141").Value = 701550.1
$ws.Range("I141").Value = 778465.25
$ws.Range("K141").Value = 2335395.75
$ws.Range("M141").Value = -2330215.75

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 309211.44
$ws.Range("I2").Value = 397384.44
$ws.Range("K2").Value = 397384.44
$ws.Range("M2").Value = -397271.44
$ws.Range("H32").Value = 3040.24
$ws.Range("I32").Value = 2825.3088
$ws.Range("J32").Value = 5128.143
$ws.Range("K32").Value = 2825.3088
$ws.Range("L32").Value = 5128.143
$ws.Range("M32").Value = -2538.3088
$ws.Range("N32").Value = -5702.143
$ws.Range("H45").Value = 1662
$ws.Range("I45").Value = 1515.4615
$ws.Range("J45").Value = 1835.1818
$ws.Range("K45").Value = 1515.4615
$ws.Range("L45").Value = 1835.1818
$ws.Range("M45").Value = -1138.4615
$ws.Range("N45").Value = -2589.1818
$ws.Range("H61").Value = 22728624
$ws.Range("I61").Value = 16667632
$ws.Range("J61").Value = 35716464
$ws.Range("K61").Value = 16667632
$ws.Range("L61").Value = 35716464
$ws.Range("M61").Value = -16667420
$ws.Range("N61").Value = -35716888
$ws.Range("H74").Value = 1110.5714
$ws.Range("I74").Value = 882.0645
$ws.Range("J74").Value = 1754.5454
$ws.Range("K74").Value = 882.0645
$ws.Range("L74").Value = 1754.5454
$ws.Range("M74").Value = -8.064499999999953
$ws.Range("N74").Value = -3502.5454
$ws.Range("H77").Value = 1110.5714
$ws.Range("I77").Value = 882.0645
$ws.Range("J77").Value = 1754.5454
$ws.Range("K77").Value = 4410.3225
$ws.Range("L77").Value = 8772.726999999999
$ws.Range("M77").Value = -42.32250000000022
$ws.Range("N77").Value = -17508.727
$ws.Range("H116").Value = 309211.44
$ws.Range("I116").Value = 397384.44
$ws.Range("K116").Value = 397384.44
$ws.Range("M116").Value = -395090.44
$ws.Range("H122").Value = 1716.4324
$ws.Range("I122").Value = 1393.0741
$ws.Range("J122").Value = 2589.5
$ws.Range("K122").Value = 4179.2223
$ws.Range("L122").Value = 7768.5
$ws.Range("M122").Value = -1729.2223
$ws.Range("N122").Value = -12668.5
$ws.Range("H130").Value = 16135.25
$ws.Range("J130").Value = 16135.25
$ws.Range("L130").Value = 16135.25
$ws.Range("N130").Value = -26175.25
$ws.Range("H136").Value = 22728624
$ws.Range("I136").Value = 16667632
$ws.Range("J136").Value = 35716464
$ws.Range("K136").Value = 50002896
$ws.Range("L136").Value = 107149392
$ws.Range("M136").Value = -50000346
$ws.Range("N136").Value = -107154492

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 309211.44
$ws.Range("I3").Value = 397384.44
$ws.Range("K3").Value = 397384.44
$ws.Range("M3").Value = -397270.44
$ws.Range("H76").Value = 75622.25
$ws.Range("J76").Value = 75622.25
$ws.Range("L76").Value = 75622.25
$ws.Range("N76").Value = -76252.25
$ws.Range("H79").Value = 75622.25
$ws.Range("J79").Value = 75622.25
$ws.Range("L79").Value = 75622.25
$ws.Range("N79").Value = -77806.25
$ws.Range("H88").Value = 12199.2
$ws.Range("J88").Value = 12749
$ws.Range("L88").Value = 12749
$ws.Range("N88").Value = -13561
$ws.Range("H91").Value = 12199.2
$ws.Range("J91").Value = 12749
$ws.Range("L91").Value = 12749
$ws.Range("N91").Value = -15557

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1624559.9
$ws.Range("I31").Value = 2101691.8
$ws.Range("J31").Value = 2311.2
$ws.Range("K31").Value = 2101691.8
$ws.Range("L31").Value = 2311.2
$ws.Range("M31").Value = -2101396.8
$ws.Range("N31").Value = -2901.2
$ws.Range("H34").Value = 1624559.9
$ws.Range("I34").Value = 2101691.8
$ws.Range("J34").Value = 2311.2
$ws.Range("K34").Value = 2101691.8
$ws.Range("L34").Value = 2311.2
$ws.Range("M34").Value = -2101489.8
$ws.Range("N34").Value = -2715.2
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H99").Value = 2500
$ws.Range("J99").Value = 2500
$ws.Range("L99").Value = 2500
$ws.Range("N99").Value = -5496
$ws.Range("H122").Value = 4616.4614
$ws.Range("I122").Value = 3996.75
$ws.Range("J122").Value = 5608
$ws.Range("K122").Value = 11990.25
$ws.Range("L122").Value = 16824
$ws.Range("M122").Value = -9540.25
$ws.Range("N122").Value = -21724
$ws.Range("H126").Value = 2500
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 1458.3658
$ws.Range("I132").Value = 1033.1072
$ws.Range("J132").Value = 2374.3076
$ws.Range("K132").Value = 3099.3216
$ws.Range("L132").Value = 7122.9228
$ws.Range("M132").Value = -569.3215999999998
$ws.Range("N132").Value = -12182.9228
$ws.Range("H134").Value = 1488.2241
$ws.Range("I134").Value = 1372.2667
$ws.Range("J134").Value = 1889.6154
$ws.Range("K134").Value = 4116.800099999999
$ws.Range("L134").Value = 5668.8462
$ws.Range("M134").Value = -1581.800099999999
$ws.Range("N134").Value = -10738.8462

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 805
$ws.Range("I5").Value = 750
$ws.Range("J5").Value = 838
$ws.Range("K5").Value = 2250
$ws.Range("L5").Value = 2514
$ws.Range("M5").Value = -2138
$ws.Range("N5").Value = -2738
$ws.Range("H122").Value = 879.1111
$ws.Range("I122").Value = 528.3077
$ws.Range("K122").Value = 4754.7693
$ws.Range("M122").Value = -2304.7693
$ws.Range("H135").Value = 805
$ws.Range("I135").Value = 750
$ws.Range("J135").Value = 838
$ws.Range("K135").Value = 6750
$ws.Range("L135").Value = 7542
$ws.Range("M135").Value = -4215
$ws.Range("N135").Value = -12612

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2399.6667
$ws.Range("I102").Value = 2392.5
$ws.Range("K102").Value = 2392.5
$ws.Range("M102").Value = -770.5
$ws.Range("H126").Value = 2177803.2
$ws.Range("I126").Value = 9262450
$ws.Range("J126").Value = 52409.2
$ws.Range("K126").Value = 27787350
$ws.Range("L126").Value = 157227.6
$ws.Range("M126").Value = -27784880
$ws.Range("N126").Value = -162167.6
$ws.Range("H132").Value = 664422.3
$ws.Range("J132").Value = 3142.4285
$ws.Range("L132").Value = 9427.2855
$ws.Range("N132").Value = -14487.2855
$ws.Range("H141").Value = 40429
$ws.Range("J141").Value = 40429
$ws.Range("L141").Value = 40429
$ws.Range("N141").Value = -50789

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3553.3845
$ws.Range("I7").Value = 3401
$ws.Range("K7").Value = 3401
$ws.Range("M7").Value = -3289
$ws.Range("H64").Value = 999999
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 999999
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H126").Value = 3553.3845
$ws.Range("I126").Value = 3401
$ws.Range("K126").Value = 10203
$ws.Range("M126").Value = -7733

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 27496.75
$ws.Range("J124").Value = 27496.75
$ws.Range("L124").Value = 27496.75
$ws.Range("N124").Value = -37316.75
$ws.Range("H126").Value = 6592.84
$ws.Range("I126").Value = 9360.385
$ws.Range("K126").Value = 28081.155
$ws.Range("M126").Value = -25611.155
$ws.Range("H132").Value = 1069.3334
$ws.Range("I132").Value = 733.5625
$ws.Range("K132").Value = 2200.6875
$ws.Range("M132").Value = 329.3125
$ws.Range("H136").Value = 11339381
$ws.Range("I136").Value = 15016435
$ws.Range("K136").Value = 45049305
$ws.Range("M136").Value = -45046755
$ws.Range("H141").Value = 72984.45
$ws.Range("J141").Value = 72984.45
$ws.Range("L141").Value = 72984.45
$ws.Range("N141").Value = -83344.45
